$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new test-case row (row 9), mirroring the existing rows' data
# (this is effectively a copy of row 8 with an updated Testcaseid and a
# few numeric-looking values typed in as real numbers instead of text).
$ws.Range("A9").Value2 = "TC008"
$ws.Range("B9").Value2 = "email1@gmail.com"
$ws.Range("C9").Value2 = "NTVNSV31"
$ws.Range("D9").Value2 = "htcglobal2019"
$ws.Range("E9").Value2 = "TN"
$ws.Range("F9").Value2 = "Bond - No Credit"
$ws.Range("G9").Value2 = 37209
$ws.Range("H9").Value2 = "Nil"
$ws.Range("I9").Value2 = "12.5% down, 11 monthly payments"
$ws.Range("J9").Value2 = "Nil"
$ws.Range("K9").Value2 = "No"
$ws.Range("L9").Value2 = "Nil"
$ws.Range("M9").Value2 = "Nil"
$ws.Range("N9").Value2 = "Nil"
$ws.Range("O9").Value2 = "Nil"
$ws.Range("P9").Value2 = "Nil"
$ws.Range("Q9").Value2 = "No"
$ws.Range("R9").Value2 = "Nil"
$ws.Range("S9").Value2 = "Nil"
$ws.Range("T9").Value2 = "Nil"
$ws.Range("U9").Value2 = "Nil"
$ws.Range("V9").Value2 = "Nil"
$ws.Range("W9").Value2 = 666195144
$ws.Range("X9").Value2 = 487956895
$ws.Range("Y9").Value2 = "Nil"
$ws.Range("Z9").Value2 = "Nil"
$ws.Range("AA9").Value2 = "Nil"
$ws.Range("AB9").Value2 = "Nil"
$ws.Range("AC9").Value2 = "Nil"
$ws.Range("AD9").Value2 = "Nil"
$ws.Range("AE9").Value2 = "Nil"
$ws.Range("AF9").Value2 = "Nil"
$ws.Range("AG9").Value2 = "Nil"

# Move / refresh the current selection, matching where the author of the
# edit ended up after typing the new row in.
$ws.Range("C15").Select() | Out-Null
